$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.490853973132971
$ws.Range("C2").Value = 7.059523881024914
$ws.Range("D2").Value = 7.23016542682854
$ws.Range("E2").Value = 11.7315790012234
$ws.Range("F2").Value = 38.72962792929081
$ws.Range("I2").Value = 30.23405911477383
$ws.Range("K2").Value = 10.65070691646888
$ws.Range("M2").Value = 14.67389947728574
$ws.Range("B3").Value = 9.349053031241512
$ws.Range("C3").Value = 6.853847722207009
$ws.Range("D3").Value = 7.215854189163946
$ws.Range("E3").Value = 11.49052634643768
$ws.Range("F3").Value = 38.26644242949627
$ws.Range("I3").Value = 30.00185048550637
$ws.Range("K3").Value = 10.55337400510517
$ws.Range("M3").Value = 14.5482192420212
$ws.Range("B4").Value = 9.26501479525829
$ws.Range("C4").Value = 6.726976198064594
$ws.Range("D4").Value = 7.206760901445559
$ws.Range("E4").Value = 11.34354466096589
$ws.Range("F4").Value = 37.98535806896161
$ws.Range("I4").Value = 29.86185222809553
$ws.Range("K4").Value = 10.49705343536169
$ws.Range("M4").Value = 14.47503259056658
$ws.Range("B5").Value = 9.231584466157724
$ws.Range("C5").Value = 6.675220826125078
$ws.Range("D5").Value = 7.202978009991683
$ws.Range("E5").Value = 11.28399334106054
$ws.Range("F5").Value = 37.87174994336257
$ws.Range("I5").Value = 29.80548439892484
$ws.Range("K5").Value = 10.47499534734592
$ws.Range("M5").Value = 14.44623943025396
$ws.Range("B6").Value = 9.226084221934904
$ws.Range("C6").Value = 6.666626419588819
$ws.Range("D6").Value = 7.202345185233777
$ws.Range("E6").Value = 11.27412827762282
$ws.Range("F6").Value = 37.85294489761796
$ws.Range("I6").Value = 29.79616676047198
$ws.Range("K6").Value = 10.47138731598167
$ws.Range("M6").Value = 14.44152144567746
$ws.Range("B7").Value = 9.264560567851172
$ws.Range("C7").Value = 6.726278300187218
$ws.Range("D7").Value = 7.206710197493325
$ws.Range("E7").Value = 11.34274001899143
$ws.Range("F7").Value = 37.98382198597133
$ws.Range("I7").Value = 29.86108922411268
$ws.Range("K7").Value = 10.49675230181827
$ws.Range("M7").Value = 14.47464006373736
$ws.Range("B8").Value = 9.441365741201324
$ws.Range("C8").Value = 6.988782632331112
$ws.Range("D8").Value = 7.225294254072631
$ws.Range("E8").Value = 11.64830309370883
$ws.Range("F8").Value = 38.56929616942924
$ws.Range("I8").Value = 30.15347924391661
$ws.Range("K8").Value = 10.61644914519881
$ws.Range("M8").Value = 14.62975710341374
$ws.Range("B9").Value = 9.809637075753002
$ws.Range("C9").Value = 7.495071363077408
$ws.Range("D9").Value = 7.25932359730022
$ws.Range("E9").Value = 12.2518895779382
$ws.Range("F9").Value = 39.73910344683858
$ws.Range("I9").Value = 30.74576070239645
$ws.Range("K9").Value = 10.87718919878496
$ws.Range("M9").Value = 14.96412290620863
$ws.Range("B10").Value = 10.08998293828789
$ws.Range("C10").Value = 7.857132262130735
$ws.Range("D10").Value = 7.282874863594933
$ws.Range("E10").Value = 12.69312999769558
$ws.Range("F10").Value = 40.60557248003055
$ws.Range("I10").Value = 31.19028275970089
$ws.Range("K10").Value = 11.082763974937
$ws.Range("M10").Value = 15.22616885684338
$ws.Range("B11").Value = 10.21892906082931
$ws.Range("C11").Value = 8.018810742887158
$ws.Range("D11").Value = 7.293278613833119
$ws.Range("E11").Value = 12.89240667693939
$ws.Range("F11").Value = 40.99995357179781
$ws.Range("I11").Value = 31.39406609418525
$ws.Range("K11").Value = 11.17890880643785
$ws.Range("M11").Value = 15.34848150138972
$ws.Range("B12").Value = 10.26790286915748
$ws.Range("C12").Value = 8.079534345704964
$ws.Range("D12").Value = 7.297174014412793
$ws.Range("E12").Value = 12.96758360663903
$ws.Range("F12").Value = 41.14921750816949
$ws.Range("I12").Value = 31.47141773517939
$ws.Range("K12").Value = 11.21565811775102
$ws.Range("M12").Value = 15.39520507833258
$ws.Range("B13").Value = 10.25734996710978
$ws.Range("C13").Value = 8.066479670605842
$ws.Range("D13").Value = 7.296337034881619
$ws.Range("E13").Value = 12.95140675591184
$ws.Range("F13").Value = 41.11707614151184
$ws.Range("I13").Value = 31.45475119735271
$ws.Range("K13").Value = 11.20772889969635
$ws.Range("M13").Value = 15.38512489131767
$ws.Range("B14").Value = 10.22295559669133
$ws.Range("C14").Value = 8.023816895085773
$ws.Range("D14").Value = 7.293599976511739
$ws.Range("E14").Value = 12.89859766165074
$ws.Range("F14").Value = 41.01223583210326
$ws.Range("I14").Value = 31.40042642911892
$ws.Range("K14").Value = 11.18192556373004
$ws.Range("M14").Value = 15.35231756547206
$ws.Range("B15").Value = 10.20190517778811
$ws.Range("C15").Value = 7.997617689673725
$ws.Range("D15").Value = 7.291917691693489
$ws.Range("E15").Value = 12.86621126508686
$ws.Range("F15").Value = 40.94800447544878
$ws.Range("I15").Value = 31.36717355539476
$ws.Range("K15").Value = 11.16616363885228
$ws.Range("M15").Value = 15.3322738549268
$ws.Range("B16").Value = 10.08157942515307
$ws.Range("C16").Value = 7.84649957935155
$ws.Range("D16").Value = 7.282188743419475
$ws.Range("E16").Value = 12.68007103168759
$ws.Range("F16").Value = 40.57979362402626
$ws.Range("I16").Value = 31.17699299670792
$ws.Range("K16").Value = 11.07653044952907
$ws.Range("M16").Value = 15.21823452057548
$ws.Range("B17").Value = 10.0080836823006
$ws.Range("C17").Value = 7.75297015379622
$ws.Range("D17").Value = 7.276141193124258
$ws.Range("E17").Value = 12.56545364517376
$ws.Range("F17").Value = 40.35388483872979
$ws.Range("I17").Value = 31.0606951230191
$ws.Range("K17").Value = 11.0221905919536
$ws.Range("M17").Value = 15.14904229988134
$ws.Range("B18").Value = 9.965947538484517
$ws.Range("C18").Value = 7.698893309955415
$ws.Range("D18").Value = 7.272633622235877
$ws.Range("E18").Value = 12.49939700971867
$ws.Range("F18").Value = 40.22397535109258
$ws.Range("I18").Value = 30.99395380603552
$ws.Range("K18").Value = 10.99118530936726
$ws.Range("M18").Value = 15.10953963589322
$ws.Range("B19").Value = 9.95170625233466
$ws.Range("C19").Value = 7.680537624500946
$ws.Range("D19").Value = 7.271440998675192
$ws.Range("E19").Value = 12.47701115411693
$ws.Range("F19").Value = 40.17999833266008
$ws.Range("I19").Value = 30.97138340833551
$ws.Range("K19").Value = 10.98073142270702
$ws.Range("M19").Value = 15.09621657444384
$ws.Range("B20").Value = 10.01589371634867
$ws.Range("C20").Value = 7.762956117759316
$ws.Range("D20").Value = 7.276787983150242
$ws.Range("E20").Value = 12.57766906765962
$ws.Range("F20").Value = 40.37793120399157
$ws.Range("I20").Value = 31.07305999525964
$ws.Range("K20").Value = 11.02794959682369
$ws.Range("M20").Value = 15.15637770186722
$ws.Range("B21").Value = 10.23305459067282
$ws.Range("C21").Value = 8.036362062403164
$ws.Range("D21").Value = 7.294405116273571
$ws.Range("E21").Value = 12.91411729429601
$ws.Range("F21").Value = 41.04303301072035
$ws.Range("I21").Value = 31.41637829306776
$ws.Range("K21").Value = 11.18949565771953
$ws.Range("M21").Value = 15.36194316578633
$ws.Range("B22").Value = 10.37579915618363
$ws.Range("C22").Value = 8.212105107874626
$ws.Range("D22").Value = 7.305660995430284
$ws.Range("E22").Value = 13.13231291977133
$ws.Range("F22").Value = 41.47720120374141
$ws.Range("I22").Value = 31.64180537875269
$ws.Range("K22").Value = 11.29704786123849
$ws.Range("M22").Value = 15.49864205475223
$ws.Range("B23").Value = 10.29955791637651
$ws.Range("C23").Value = 8.118597138785216
$ws.Range("D23").Value = 7.299677030823412
$ws.Range("E23").Value = 13.01603697400823
$ws.Range("F23").Value = 41.24556077354951
$ws.Range("I23").Value = 31.52140875375381
$ws.Range("K23").Value = 11.23947688772011
$ws.Range("M23").Value = 15.42548170633237
$ws.Range("B24").Value = 10.01236242934037
$ws.Range("C24").Value = 7.758442410078962
$ws.Range("D24").Value = 7.276495665183575
$ws.Range("E24").Value = 12.57214697257125
$ws.Range("F24").Value = 40.3670599318532
$ws.Range("I24").Value = 31.06746945983827
$ws.Range("K24").Value = 11.02534521462411
$ws.Range("M24").Value = 15.15306050118375
$ws.Range("B25").Value = 9.708048747748441
$ws.Range("C25").Value = 7.359527348687117
$ws.Range("D25").Value = 7.250373940265028
$ws.Range("E25").Value = 12.08864971251944
$ws.Range("F25").Value = 39.42095730829814
$ws.Range("I25").Value = 30.58372714413349
$ws.Range("K25").Value = 10.8040528742338
$ws.Range("M25").Value = 14.87063865240753
